# Mark completed objectives so far with an asterisk in column D.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rows 2, 3, 4 (first three CORE FEATURES items) and row 47
# (Manage your project in a private GIT repo) are flagged as completed.
$ws.Range("D2").Value = "*"
$ws.Range("D3").Value = "*"
$ws.Range("D4").Value = "*"
$ws.Range("D47").Value = "*"

# Scroll the view down to where editing left off, and leave the
# selection on D46 to match where the author was working.
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D46").Select()
